$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.328558632615739
$ws.Range("C3").Value = -1.435981453719049
$ws.Range("E3").Value = -0.7704417043119083
$ws.Range("C4").Value = 0.406633294022174
$ws.Range("E4").Value = -0.1034614224434405
$ws.Range("C5").Value = 2.185496833134781
$ws.Range("E5").Value = 0.6652762968575532
$ws.Range("C6").Value = 0.8574941660507873
$ws.Range("E6").Value = 1.693469135756587
$ws.Range("C7").Value = 0.1494732105682406
$ws.Range("E7").Value = 0.8024032015999882
$ws.Range("C8").Value = 1.656936590801972
$ws.Range("E8").Value = 0.922773818606859
$ws.Range("C9").Value = 1.346932828201242
$ws.Range("E9").Value = 1.364302026343633
$ws.Range("C10").Value = 1.745747589686109
$ws.Range("E10").Value = 1.644798626926303
$ws.Range("C11").Value = 1.843649045891893
$ws.Range("E11").Value = 1.741128155516525
$ws.Range("C12").Value = 2.463589365374652
$ws.Range("E12").Value = 2.149194501693219
$ws.Range("C13").Value = 1.332860091726285
$ws.Range("E13").Value = 1.799885362733189
$ws.Range("C14").Value = 0.2336391425753925
$ws.Range("E14").Value = 0.9207450904090253
$ws.Range("C15").Value = -2.010709456685855
$ws.Range("E15").Value = -1.14257141002756
$ws.Range("C16").Value = 1.5286818008164
$ws.Range("E16").Value = -0.7118141543333012
$ws.Range("C17").Value = 0.001079933351455509
$ws.Range("E17").Value = 0.6889047703476203
$ws.Range("C18").Value = -0.00209793826797533
$ws.Range("E18").Value = 0.286657616500996
$ws.Range("C19").Value = 0.501314651583451
$ws.Range("E19").Value = 0.2155158706220295
